$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G (header "K") for rows 2-14
$values = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 3
    6  = 2
    7  = 0
    8  = 0
    9  = 1
    10 = 3
    11 = 2
    12 = 0
    13 = 0
    14 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
